$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of columns H and I for rows 1 through 17
# (header row 1 plus the 16 data rows), matching the reordered
# IsSourceOf / IsDerivedFrom columns in the plot test expected answers.
for ($r = 1; $r -le 17; $r++) {
    $hCell = $ws.Cells.Item($r, 8)
    $iCell = $ws.Cells.Item($r, 9)

    $hVal = $hCell.Value2
    $iVal = $iCell.Value2

    $hCell.Value = $iVal
    $iCell.Value = $hVal
}

# Update the saved selection to match the new active cell.
$ws.Range("K10").Select()
